# edit.ps1 - applies the "LLM_Airpod Max.docx" edit:
#   1) Drops the parenthetical " (avg 3.3/5; NPS -12)" from the Summary
#      insight bullet (keeping the trailing period).
#   2/3/4) Cleans up the grammar-checker run-splits (and removes the
#      now-stale <w:proofErr/> markers) in the Positive/Neutral/Negative
#      reviews(%) bullets by re-typing each line as a single run.
#   5) Same clean-up for the "Positive: Many praise ..." bullet.
#   6) Same clean-up for the "Heavy fit; pressure ..." bullet.
#   7) Same clean-up for the long "Prioritize reliability fixes ..."
#      recommendation paragraph.

$d = $word.ActiveDocument
$minus = [char]0x2212

# ---------------------------------------------------------------------
# 1) Summary insight bullet: remove " (avg 3.3/5; NPS -12)" before the
#    final period.
# ---------------------------------------------------------------------
$oldSummary = "Summary insight: AirPods Max delight many with rich sound and ANC but reliability and comfort complaints depress satisfaction (avg 3.3/5; NPS ${minus}12)."
$newSummary = "Summary insight: AirPods Max delight many with rich sound and ANC but reliability and comfort complaints depress satisfaction."
$d.Content.Find.Execute($oldSummary, $true, $false, $false, $false, $false, $true, 1, $false, $newSummary, 2) | Out-Null

# ---------------------------------------------------------------------
# 2/3/4) Positive / Neutral / Negative reviews(%) bullets.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Positive reviews(%): 36%", $true, $false, $false, $false, $false, $true, 1, $false, "Positive reviews(%): 36%", 2) | Out-Null
$d.Content.Find.Execute("Neutral reviews(%): 24%", $true, $false, $false, $false, $false, $true, 1, $false, "Neutral reviews(%): 24%", 2) | Out-Null
$d.Content.Find.Execute("Negative reviews(%): 40%", $true, $false, $false, $false, $false, $true, 1, $false, "Negative reviews(%): 40%", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) "Positive: Many praise clear, rich sound ..." bullet.
# ---------------------------------------------------------------------
$posMany = "Positive: Many praise clear, rich sound with good detail and an engaging presentation."
$d.Content.Find.Execute($posMany, $true, $false, $false, $false, $false, $true, 1, $false, $posMany, 2) | Out-Null

# ---------------------------------------------------------------------
# 6) "Heavy fit; pressure causes discomfort ..." bullet.
# ---------------------------------------------------------------------
$heavyFit = "Heavy fit; pressure causes discomfort for some users."
$d.Content.Find.Execute($heavyFit, $true, $false, $false, $false, $false, $true, 1, $false, $heavyFit, 2) | Out-Null

# ---------------------------------------------------------------------
# 7) "Prioritize reliability fixes ..." recommendation paragraph.
# ---------------------------------------------------------------------
$prioritize = "Prioritize reliability fixes: strengthen Bluetooth stack and handoff logic, add robust recovery without hard resets, and improve charging detection circuitry. Reduce weight and clamp force by revisiting headband tension and cup materials, and offer thicker/softer pads. Enhance ANC consistency via firmware tuning across modes. Add clearer on-device status/alerts for charge and connection. Finally, expand moisture ingress protection around drivers and connectors to mitigate long-term failures."
$d.Content.Find.Execute($prioritize, $true, $false, $false, $false, $false, $true, 1, $false, $prioritize, 2) | Out-Null
